$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellref, $val) {
    # Force the cell to be treated as text so Excel does not
    # reinterpret number-like strings (preserves exact formatting,
    # e.g. trailing zeros like "1.380") - then restore default style.
    $ws.Range($cellref).NumberFormat = "@"
    $ws.Range($cellref).Value = $val
    $ws.Range($cellref).Style = "Normal"
}

$ws.Range("D2").Value = "24.406.53"
$ws.Range("E2").Value = "  -1.55%  "

$ws.Range("D3").Value = "1.655.13"
$ws.Range("E3").Value = "  -2.63%  "

Set-TextValue "D4" "1.001"
$ws.Range("E4").Value = "  +0.08%  "

Set-TextValue "D5" "311.29"
$ws.Range("E5").Value = "  -1.12%  "

$ws.Range("E6").Value = "  +0.07%  "

Set-TextValue "D7" "0.3916"
$ws.Range("E7").Value = "  -1.65%  "

Set-TextValue "D8" "0.3916"
$ws.Range("E8").Value = "  -2.88%  "

Set-TextValue "D9" "1.002"
$ws.Range("E9").Value = "  -0.08%  "

Set-TextValue "D10" "1.380"
$ws.Range("E10").Value = "  -6.12%  "

Set-TextValue "D11" "49.95"
$ws.Range("E11").Value = "  -6.39%  "

Set-TextValue "D12" "0.08542"
$ws.Range("E12").Value = "  -2.89%  "

Set-TextValue "D13" "24.83"
$ws.Range("E13").Value = "  -4.95%  "

Set-TextValue "D14" "7.240"
$ws.Range("E14").Value = "  -3.86%  "

$ws.Range("E15").Value = "  -3.12%  "

Set-TextValue "D16" "7.614"
$ws.Range("E16").Value = "  -4.67%  "

$ws.Range("D17").Value = "1.654.47"
$ws.Range("E17").Value = "  -2.23%  "

Set-TextValue "D18" "93.53"
$ws.Range("E18").Value = "  -2.28%  "

Set-TextValue "D19" "0.06948"
$ws.Range("E19").Value = "  -3.21%  "

Set-TextValue "D20" "20.88"
$ws.Range("E20").Value = "  +0.09%  "

Set-TextValue "D21" "7.017"
$ws.Range("E21").Value = "  -4.44%  "

$ws.Range("E22").Value = "  +0.10%  "

$ws.Range("E23").Value = "  -3.92%  "

$ws.Range("D24").Value = "24.398.96"

Set-TextValue "D25" "2.338"
$ws.Range("E25").Value = "  -1.49%  "

Set-TextValue "D26" "2.776"
$ws.Range("E26").Value = "  -5.20%  "

Set-TextValue "D27" "22.71"
$ws.Range("E27").Value = "  -1.90%  "

Set-TextValue "D28" "159.36"
$ws.Range("E28").Value = "  -1.43%  "

Set-TextValue "D29" "5.686"
$ws.Range("E29").Value = "  -7.81%  "

Set-TextValue "D30" "145.32"
$ws.Range("E30").Value = "  +0.97%  "

Set-TextValue "D31" "8.134"
$ws.Range("E31").Value = "  -4.02%  "

Set-TextValue "D32" "2.624"
$ws.Range("E32").Value = "  +9.05%  "

$ws.Range("D33").Value = "1.841.56"
$ws.Range("E33").Value = "  -3.93%  "

Set-TextValue "D34" "1.017"
$ws.Range("E34").Value = "  -1.98%  "

Set-TextValue "D35" "0.08184"
$ws.Range("E35").Value = "  -5.22%  "

Set-TextValue "D36" "0.03014"
$ws.Range("E36").Value = "  -4.94%  "

Set-TextValue "D37" "6.863"
$ws.Range("E37").Value = "  -6.31%  "

Set-TextValue "D38" "0.2769"
$ws.Range("E38").Value = "  -2.42%  "

Set-TextValue "D39" "0.09455"
$ws.Range("E39").Value = "  +0.07%  "

Set-TextValue "D40" "10.24"
$ws.Range("E40").Value = "  -5.07%  "

Set-TextValue "D41" "1.489"
$ws.Range("E41").Value = "  +0.76%  "

Set-TextValue "D42" "0.7821"
$ws.Range("E42").Value = "  -5.77%  "

Set-TextValue "D43" "13.44"
$ws.Range("E43").Value = "  -5.44%  "

Set-TextValue "D44" "16.42"
$ws.Range("E44").Value = "  -7.62%  "

Set-TextValue "D45" "2.559"
$ws.Range("E45").Value = "  -5.25%  "

Set-TextValue "D46" "0.7042"
$ws.Range("E46").Value = "  -5.19%  "

Set-TextValue "D47" "4.150"

Set-TextValue "D48" "0.08630"
$ws.Range("E48").Value = "  +3.08%  "

$ws.Range("E49").Value = "  +0.03%  "

Set-TextValue "D50" "1.310"
$ws.Range("E50").Value = "  -6.28%  "

Set-TextValue "D51" "136.89"
$ws.Range("E51").Value = "  -2.11%  "
